# Update last 30 days report: remove the "Oscilar / Enterprise AE x5 /
# JOAQUIN SALCEDO" entry (row 5). Deleting the entire row shifts all
# subsequent rows up by one, which matches the new data shown in the diff
# and shrinks the used range from A1:F16 down to A1:F15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(5).Delete()
